# NIT-9011856825.xlsx — "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
#
# The account-statement table (rows 16-42) is rebuilt:
#   - PEDRO ANTONIO MARMOL MARMOL (CC 73157434): periods 2309-2506 -> 2309-2507 (23 periods)
#   - LUIS DAVID LLAMAS RUIZ (CC 1047471479): periods 2309-2506 -> 2309-2401 only (5 periods)
#   - CARENTH ROMERO HIDALGO (CC 67040568): brand-new worker, periods 2504-2505 (2 periods)
# Table grows from 27 data rows to 30, so 3 rows are inserted before the old
# closing row (old row 42) and the whole grid (incl. the footer block) is
# rewritten/shifted accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room: insert 3 rows before the old "closing" row (42) ---------
$ws.Range("42:44").Insert()

# Copy the regular data-row formatting (row 41, untouched by the insert) onto
# the 3 freshly inserted rows so borders/fonts/number-formats match the rest
# of the table instead of Excel's generic insert-default style.
$ws.Range("B41:J41").Copy()
$ws.Range("B42:J44").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 2. Rewrite every data row (16-45) ---------------------------------------
$data = @(
    @(16, "CC", "73157434", "PEDRO ANTONIO MARMOL MARMOL", "2507", 46400, 1160000),
    @(17, "CC", "73157434", "PEDRO ANTONIO MARMOL MARMOL", "2506", 46400, 1160000),
    @(18, "CC", "73157434", "PEDRO ANTONIO MARMOL MARMOL", "2505", 46400, 1160000),
    @(19, "CC", "73157434", "PEDRO ANTONIO MARMOL MARMOL", "2504", 46400, 1160000),
    @(20, "CC", "73157434", "PEDRO ANTONIO MARMOL MARMOL", "2503", 46400, 1160000),
    @(21, "CC", "73157434", "PEDRO ANTONIO MARMOL MARMOL", "2502", 46400, 1160000),
    @(22, "CC", "73157434", "PEDRO ANTONIO MARMOL MARMOL", "2501", 46400, 1160000),
    @(23, "CC", "73157434", "PEDRO ANTONIO MARMOL MARMOL", "2412", 46400, 1160000),
    @(24, "CC", "73157434", "PEDRO ANTONIO MARMOL MARMOL", "2411", 46400, 1160000),
    @(25, "CC", "73157434", "PEDRO ANTONIO MARMOL MARMOL", "2410", 46400, 1160000),
    @(26, "CC", "73157434", "PEDRO ANTONIO MARMOL MARMOL", "2409", 46400, 1160000),
    @(27, "CC", "73157434", "PEDRO ANTONIO MARMOL MARMOL", "2408", 46400, 1160000),
    @(28, "CC", "73157434", "PEDRO ANTONIO MARMOL MARMOL", "2407", 46400, 1160000),
    @(29, "CC", "73157434", "PEDRO ANTONIO MARMOL MARMOL", "2406", 46400, 1160000),
    @(30, "CC", "73157434", "PEDRO ANTONIO MARMOL MARMOL", "2405", 46400, 1160000),
    @(31, "CC", "73157434", "PEDRO ANTONIO MARMOL MARMOL", "2404", 46400, 1160000),
    @(32, "CC", "73157434", "PEDRO ANTONIO MARMOL MARMOL", "2403", 46400, 1160000),
    @(33, "CC", "73157434", "PEDRO ANTONIO MARMOL MARMOL", "2402", 46400, 1160000),
    @(34, "CC", "73157434", "PEDRO ANTONIO MARMOL MARMOL", "2401", 46400, 1160000),
    @(35, "CC", "73157434", "PEDRO ANTONIO MARMOL MARMOL", "2312", 46400, 1160000),
    @(36, "CC", "73157434", "PEDRO ANTONIO MARMOL MARMOL", "2311", 46400, 1160000),
    @(37, "CC", "73157434", "PEDRO ANTONIO MARMOL MARMOL", "2310", 46400, 1160000),
    @(38, "CC", "73157434", "PEDRO ANTONIO MARMOL MARMOL", "2309", 1547, 1160000),
    @(39, "CC", "1047471479", "LUIS DAVID LLAMAS RUIZ", "2401", 44854, 1160000),
    @(40, "CC", "1047471479", "LUIS DAVID LLAMAS RUIZ", "2312", 46400, 1160000),
    @(41, "CC", "1047471479", "LUIS DAVID LLAMAS RUIZ", "2311", 46400, 1160000),
    @(42, "CC", "1047471479", "LUIS DAVID LLAMAS RUIZ", "2310", 46400, 1160000),
    @(43, "CC", "1047471479", "LUIS DAVID LLAMAS RUIZ", "2309", 1547, 1160000),
    @(44, "CC", "67040568", "CARENTH ROMERO HIDALGO", "2505", 200000, 5000000),
    @(45, "CC", "67040568", "CARENTH ROMERO HIDALGO", "2504", 200000, 5000000)
)

$totalMora = 0
foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("B$r").Value2 = $row[1]
    $ws.Range("C$r").Value2 = $row[2]
    $ws.Range("D$r").Value2 = $row[3]
    $ws.Range("E$r").Value2 = $row[4]
    $ws.Range("F$r").Value2 = $row[5]
    $ws.Range("G$r").Value2 = $row[6]
    $totalMora = $totalMora + $row[5]
}

# --- 3. Summary header ------------------------------------------------------
$ws.Range("E11").Value2 = $totalMora        # VALOR MORA total (1547+...+200000 = 1607948)
$ws.Range("C13").Value2 = 3                 # Cant. Trabajadores (2 -> 3)
$ws.Range("F13").Value2 = 23                 # Cant. Periodos, distinct periods (22 -> 23)
